# Fruta / hortaliza, semanal
# Insert a new weekly price record as row 232, pushing all subsequent rows
# (the former rows 232-282) down by one to become rows 233-283.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 232 (shifts rows 232..282 -> 233..283)
$ws.Rows.Item(232).Insert()

# Populate the newly inserted row 232 with the new weekly record
$ws.Cells.Item(232, 1).Value = 7
$ws.Cells.Item(232, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(232, 3).Value = "Ñuble"
$ws.Cells.Item(232, 4).Value = 45258
$ws.Cells.Item(232, 5).Value = 16
$ws.Cells.Item(232, 6).Value = 100112028
$ws.Cells.Item(232, 7).Value = "Sandia"
$ws.Cells.Item(232, 8).Value = "Sin especificar"
$ws.Cells.Item(232, 9).Value = "Primera"
$ws.Cells.Item(232, 10).Value = 300
$ws.Cells.Item(232, 11).Value = 900
$ws.Cells.Item(232, 12).Value = 900
$ws.Cells.Item(232, 13).Value = 900
$ws.Cells.Item(232, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(232, 15).Value = "Perú"
$ws.Cells.Item(232, 16).Value = 900
$ws.Cells.Item(232, 17).Value = 1
$ws.Cells.Item(232, 18).Value = "Hortaliza"
